# Records.xlsx edit: rename sheet, update headers, resize columns, update active sheet/selection.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsRules = $wb.Worksheets.Item(2)

# --- Rename "Benchmark Sequences" -> "Benchmark_Sequences" (also updates the
#     defined name formula automatically since the sheet name no longer needs quoting).
$wsData.Name = "Benchmark_Sequences"

# --- Header text tweaks on the data sheet (row 2).
$wsData.Range("E2").Value = "Scene"
$wsData.Range("F2").Value = "Duration"
$wsData.Range("G2").Value = "Motion_Dyn"

# --- Column width changes (values are in "characters"; COM rounds to the
#     nearest pixel internally, same as Excel does).
$wsData.Columns.Item(1).ColumnWidth = 14.666666666666666
$wsData.Columns.Item(2).ColumnWidth = 25.333333333333332
$wsData.Columns.Item(3).ColumnWidth = 45.666666666666664
$wsData.Range($wsData.Columns.Item(4), $wsData.Columns.Item(10)).ColumnWidth = 25.333333333333332
$wsData.Range($wsData.Columns.Item(14), $wsData.Columns.Item(1023)).ColumnWidth = 20.0
$wsData.Range($wsData.Columns.Item(1024), $wsData.Columns.Item(1025)).ColumnWidth = 8.833333333333334

$wsRules.Columns.Item(1).ColumnWidth = 20.166666666666668
$wsRules.Range($wsRules.Columns.Item(2), $wsRules.Columns.Item(4)).ColumnWidth = 50.666666666666664
$wsRules.Range($wsRules.Columns.Item(5), $wsRules.Columns.Item(1025)).ColumnWidth = 13.333333333333334

# --- Active sheet / selection.
$null = $wsData.Activate()
$null = $wsData.Range("J2").Select()
